# DiarioTabela.xlsx - "Criação Home Page Part-1" commit
#
# - Adds a "Grupo 10" title in A2 (above the existing Data/Descricao/Tarefa
#   table that lives in rows 4-6).
# - Appends a new diary entry in row 7:
#     A7 = 29/10/2017
#     B7 = a free-text note that starts with "-" so Excel treats it as a
#          formula; since "Definicao" etc. aren't valid names/functions it
#          evaluates to a #NAME? error (exactly as in the source workbook).
#     C7 = "Criação Home Page"
#     D7 = "2h" (re-using the same time-spent text already used in D6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Grupo 10" heading above the table.
$ws.Range("A2").Value = "Grupo 10"

# New diary row.
$ws.Range("A7").Value = "29/10/2017"
$ws.Range("B7").Formula = "=-Definicao do Template de todas as Páginas do Site, United Template     -Definicão dos Atributos da Página Home -Criação de algumas Views"
$ws.Range("C7").Value = "Criação Home Page"
$ws.Range("D7").Value = "2h"

# Match the author's final selection/cursor position.
$ws.Range("D7").Select() | Out-Null
